$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; this pushes every existing row down by one.
$ws.Rows.Item(1).Insert()

# The floating screenshot picture is anchored "move but don't size with cells", so it
# does not automatically follow the row insert in this engine - shift it down by
# exactly one default row height (14.4 pt) so its anchor row advances by one, just
# like the rest of the sheet did.
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top + 14.4

# ---------------------------------------------------------------------------
# Row 1 (new): "Chromium | Web" sub-header, keeping the old xpath selector columns
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Chromium"
$ws.Range("B1").Value2 = "Web"
$ws.Range("C1").Value2 = "xpath=//*/li[@class = 'nav-item']/a[text()=' Profil']"
$ws.Range("D1").Value2 = "xpath=//*/li[@class = 'nav-item']/a[text()=' Abwesenheiten']"
$ws.Range("E1").Value2 = "xpath=//*/li[@class = 'nav-item']/a[text()=' Benachrichtigungen']"
$ws.Range("F1").Value2 = "Action"

# ---------------------------------------------------------------------------
# Row 2 (new): "Pixel9Pro_API35 | Mobile" sub-header, with empty highlighted cells
# ---------------------------------------------------------------------------
$ws.Range("A2").Value2 = "Pixel9Pro_API35"
$ws.Range("B2").Value2 = "Mobile"
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value2 = "Action"

# ---------------------------------------------------------------------------
# Row 3: main column header row (used to be row 1), now bold with a gray fill
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = "Record / Control "
$ws.Range("B3").Value2 = "Modus"
$ws.Range("C3").Value2 = "pagProfil"
$ws.Range("D3").Value2 = "pagAbwesenheiten"
$ws.Range("E3").Value2 = "pagBenachrichtigungen"
$ws.Range("F3").Value2 = "Action"

# --- "Text" number format for the sub-header labels in A1:B2 ---
$ws.Range("A1:B1").NumberFormat = "@"
$ws.Range("A2:B2").NumberFormat = "@"

# --- Yellow highlight fill for the still-empty C2:E2 cells ---
$ws.Range("C2:E2").Interior.Color = 65535

# --- Bold + gray header style for row 3. Build it once on a scratch cell
#     (seeded with the same yellow fill already in use) and then copy/paste the
#     resolved format onto the real range, so the engine only ever resolves the
#     single, final combined style instead of leaving unused in-between states
#     behind in the style table. ---
$scratch = $ws.Range("Z100")
$scratch.Interior.Color = 65535
$scratch.NumberFormat = "@"
$scratch.Font.Bold = $true
$scratch.Interior.ThemeColor = 4
$scratch.Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = 0

# Leave the final selection on D15, matching where editing ended up.
$ws.Range("D15").Select()
